$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 578, pushing all existing
# data (old rows 578-635) down to rows 580-637.
$ws.Rows.Item(578).Insert()
$ws.Rows.Item(578).Insert()

# --- New row 578: Hass / Primera, new price observation ---
$ws.Cells.Item(578, 1).Value = 11
$ws.Cells.Item(578, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(578, 3).Value = "Bíobío"
$ws.Cells.Item(578, 4).Value = 44783
$ws.Cells.Item(578, 5).Value = 8
$ws.Cells.Item(578, 6).Value = "Fruta"
$ws.Cells.Item(578, 7).Value = 100106
$ws.Cells.Item(578, 8).Value = "Oleaginosos"
$ws.Cells.Item(578, 9).Value = 100106002
$ws.Cells.Item(578, 10).Value = "Palta"
$ws.Cells.Item(578, 11).Value = "Hass"
$ws.Cells.Item(578, 12).Value = "Primera"
$ws.Cells.Item(578, 13).Value = 250
$ws.Cells.Item(578, 14).Value = 2600
$ws.Cells.Item(578, 15).Value = 2700
$ws.Cells.Item(578, 16).Value = 2660
$ws.Cells.Item(578, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(578, 18).Value = "Perú"
$ws.Cells.Item(578, 19).Value = 2660
$ws.Cells.Item(578, 20).Value = 1

# --- New row 579: Hass / Segunda, new price observation ---
$ws.Cells.Item(579, 1).Value = 11
$ws.Cells.Item(579, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(579, 3).Value = "Bíobío"
$ws.Cells.Item(579, 4).Value = 44783
$ws.Cells.Item(579, 5).Value = 8
$ws.Cells.Item(579, 6).Value = "Fruta"
$ws.Cells.Item(579, 7).Value = 100106
$ws.Cells.Item(579, 8).Value = "Oleaginosos"
$ws.Cells.Item(579, 9).Value = 100106002
$ws.Cells.Item(579, 10).Value = "Palta"
$ws.Cells.Item(579, 11).Value = "Hass"
$ws.Cells.Item(579, 12).Value = "Segunda"
$ws.Cells.Item(579, 13).Value = 250
$ws.Cells.Item(579, 14).Value = 2300
$ws.Cells.Item(579, 15).Value = 2300
$ws.Cells.Item(579, 16).Value = 2300
$ws.Cells.Item(579, 17).Value = "$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(579, 18).Value = "Perú"
$ws.Cells.Item(579, 19).Value = 2300
$ws.Cells.Item(579, 20).Value = 1
